$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (login/logout/working time now stored as quoted text) ---
$ws.Range("D2").Value = '"2:35PM"'
$ws.Range("E2").Value = '"6:00PM"'
$ws.Range("F2").Value = '"03:30"'

# --- Insert new attendance record in row 3 ---
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Avijit"
$ws.Range("C3").Value = '"29/10/21"'
$ws.Range("D3").Value = '"14:14PM"'
$ws.Range("E3").Value = '"18:55PM"'
$ws.Range("F3").Value = '"04:30"'

# --- Number format for the workingtime column (h:mm, no seconds/AM-PM) ---
$ws.Range("F2").NumberFormat = "h:mm"
$ws.Range("F3").NumberFormat = "h:mm"

# --- Widen the Date column to fit the longer "dd/mm/yy" text ---
$ws.Columns("C").ColumnWidth = 9.6

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Final selection matches the authored state ---
$ws.Range("C4").Select() | Out-Null
